# Updates the cryptos list: refreshed prices and 1h volume percentages,
# plus a re-ordering of a couple of rows (Litecoin/Dai swap positions 22/23,
# and USDe/Aptos swap positions 32/33) to reflect the latest ranking snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D:E) hold free-form text values (e.g. "68.167.13",
# "0.603", "  +0.14%  ") that must remain text rather than being
# auto-converted to numbers/dates by Excel's type inference. Force the
# relevant range to Text format before writing, then restore the original
# (unstyled) Normal style once all values have been written.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$updates = @(
    @{ Cell = 'D2'; Value = '68.053.20' },
    @{ Cell = 'E2'; Value = '  +0.23%  ' },
    @{ Cell = 'D3'; Value = '3.269.09' },
    @{ Cell = 'E3'; Value = '  +0.88%  ' },
    @{ Cell = 'E4'; Value = '  -0.14%  ' },
    @{ Cell = 'D5'; Value = '583.80' },
    @{ Cell = 'E5'; Value = '  +0.30%  ' },
    @{ Cell = 'D6'; Value = '184.59' },
    @{ Cell = 'E6'; Value = '  +2.63%  ' },
    @{ Cell = 'E7'; Value = '  +0.02%  ' },
    @{ Cell = 'D8'; Value = '0.602' },
    @{ Cell = 'E8'; Value = '  +1.45%  ' },
    @{ Cell = 'E9'; Value = '  -2.79%  ' },
    @{ Cell = 'E10'; Value = '  -0.02%  ' },
    @{ Cell = 'D11'; Value = '0.409' },
    @{ Cell = 'E11'; Value = '  -2.53%  ' },
    @{ Cell = 'D12'; Value = '3.838.66' },
    @{ Cell = 'E12'; Value = '  +0.42%  ' },
    @{ Cell = 'E13'; Value = '  +1.13%  ' },
    @{ Cell = 'D14'; Value = '27.49' },
    @{ Cell = 'E14'; Value = '  -1.82%  ' },
    @{ Cell = 'D15'; Value = '68.071.94' },
    @{ Cell = 'E15'; Value = '  +0.19%  ' },
    @{ Cell = 'E16'; Value = '  -0.99%  ' },
    @{ Cell = 'D17'; Value = '3.234.57' },
    @{ Cell = 'E17'; Value = '  -1.27%  ' },
    @{ Cell = 'D18'; Value = '5.74' },
    @{ Cell = 'E18'; Value = '  -0.87%  ' },
    @{ Cell = 'D19'; Value = '13.33' },
    @{ Cell = 'E19'; Value = '  -0.84%  ' },
    @{ Cell = 'D20'; Value = '415.76' },
    @{ Cell = 'E20'; Value = '  +6.26%  ' },
    @{ Cell = 'D21'; Value = '7.54' },
    @{ Cell = 'E21'; Value = '  -1.25%  ' },
    @{ Cell = 'B22'; Value = 'Litecoin' },
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' },
    @{ Cell = 'D22'; Value = '71.48' },
    @{ Cell = 'E22'; Value = '  +0.05%  ' },
    @{ Cell = 'B23'; Value = 'Dai' },
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Cell = 'D23'; Value = '0.999' },
    @{ Cell = 'E23'; Value = '  -0.06%  ' },
    @{ Cell = 'D24'; Value = '0.508' },
    @{ Cell = 'E24'; Value = '  -1.07%  ' },
    @{ Cell = 'D25'; Value = '0.0000117' },
    @{ Cell = 'E25'; Value = '  -0.38%  ' },
    @{ Cell = 'E26'; Value = '  -0.51%  ' },
    @{ Cell = 'E27'; Value = '  -0.62%  ' },
    @{ Cell = 'D28'; Value = '0.991' },
    @{ Cell = 'E28'; Value = '  -0.76%  ' },
    @{ Cell = 'D29'; Value = '1.96' },
    @{ Cell = 'E29'; Value = '  -1.06%  ' },
    @{ Cell = 'D30'; Value = '22.69' },
    @{ Cell = 'E30'; Value = '  -0.78%  ' },
    @{ Cell = 'D31'; Value = '5.46' },
    @{ Cell = 'E31'; Value = '  -2.74%  ' },
    @{ Cell = 'B32'; Value = 'USDe' },
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde' },
    @{ Cell = 'D32'; Value = '0.999' },
    @{ Cell = 'E32'; Value = '  +0.01%  ' },
    @{ Cell = 'B33'; Value = 'Aptos' },
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' },
    @{ Cell = 'D33'; Value = '6.88' },
    @{ Cell = 'E33'; Value = '  -2.74%  ' },
    @{ Cell = 'E34'; Value = '  -1.37%  ' },
    @{ Cell = 'D35'; Value = '163.59' },
    @{ Cell = 'E35'; Value = '  -0.44%  ' },
    @{ Cell = 'E36'; Value = '  -2.05%  ' },
    @{ Cell = 'E37'; Value = '  -0.17%  ' },
    @{ Cell = 'D38'; Value = '27.09' },
    @{ Cell = 'E38'; Value = '  +3.11%  ' },
    @{ Cell = 'D39'; Value = '0.798' },
    @{ Cell = 'E39'; Value = '  -2.75%  ' },
    @{ Cell = 'D40'; Value = '4.47' },
    @{ Cell = 'E40'; Value = '  -2.32%  ' },
    @{ Cell = 'D41'; Value = '6.30' },
    @{ Cell = 'E41'; Value = '  -2.94%  ' },
    @{ Cell = 'D42'; Value = '2.649.58' },
    @{ Cell = 'E42'; Value = '  +2.23%  ' },
    @{ Cell = 'D43'; Value = '40.77' },
    @{ Cell = 'E43'; Value = '  -1.41%  ' },
    @{ Cell = 'D44'; Value = '0.0677' },
    @{ Cell = 'E44'; Value = '  -0.79%  ' },
    @{ Cell = 'D45'; Value = '2.43' },
    @{ Cell = 'E45'; Value = '  -0.30%  ' },
    @{ Cell = 'D46'; Value = '337.02' },
    @{ Cell = 'E46'; Value = '  -0.60%  ' },
    @{ Cell = 'D47'; Value = '24.46' },
    @{ Cell = 'E47'; Value = '  -0.17%  ' },
    @{ Cell = 'E48'; Value = '  -2.02%  ' },
    @{ Cell = 'D49'; Value = '6.28' },
    @{ Cell = 'E49'; Value = '  +0.01%  ' },
    @{ Cell = 'D50'; Value = '0.978' },
    @{ Cell = 'E50'; Value = '  +0.71%  ' },
    @{ Cell = 'E51'; Value = '  -1.21%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

$dataRange.Style = "Normal"
